$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.63600841952526
$ws.Range("C2").Value = 4.562363476012454
$ws.Range("D2").Value = 8.577642651360105
$ws.Range("E2").Value = 10.07149674727938
$ws.Range("F2").Value = 43.39203069139788
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 35.38290916615529
$ws.Range("K2").Value = 12.91868258295409
$ws.Range("L2").Value = 10.38218828623071

$ws.Range("B3").Value = 14.53035442462118
$ws.Range("C3").Value = 4.327100145637687
$ws.Range("D3").Value = 8.569480945865621
$ws.Range("E3").Value = 10.07816099576454
$ws.Range("F3").Value = 43.04555017766453
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 35.21440922354903
$ws.Range("K3").Value = 12.84548908304078
$ws.Range("L3").Value = 10.37515367810474

$ws.Range("B4").Value = 14.47017483856947
$ws.Range("C4").Value = 4.174976701742558
$ws.Range("D4").Value = 8.564290945350967
$ws.Range("E4").Value = 10.08373185145541
$ws.Range("F4").Value = 42.83813177391556
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 35.11473204334627
$ws.Range("K4").Value = 12.80449573145671
$ws.Range("L4").Value = 10.37281800431324

$ws.Range("B5").Value = 14.4468583345201
$ws.Range("C5").Value = 4.111083119008512
$ws.Range("D5").Value = 8.562130747608352
$ws.Range("E5").Value = 10.08637415429253
$ws.Range("F5").Value = 42.75500481761795
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 35.07508235710151
$ws.Range("K5").Value = 12.78879970875627
$ws.Range("L5").Value = 10.37236603812561

$ws.Range("B6").Value = 14.44306030516678
$ws.Range("C6").Value = 4.100359658481573
$ws.Range("D6").Value = 8.561769307758524
$ws.Range("E6").Value = 10.0868353877573
$ws.Range("F6").Value = 42.74128775210465
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 35.06855763001283
$ws.Range("K6").Value = 12.78625478861131
$ws.Range("L6").Value = 10.37232120365971

$ws.Range("B7").Value = 14.46985546211245
$ws.Range("C7").Value = 4.174122672118666
$ws.Range("D7").Value = 8.564261995716352
$ws.Range("E7").Value = 10.08376597947069
$ws.Range("F7").Value = 42.83700495902758
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 35.11419336476484
$ws.Range("K7").Value = 12.80427994323925
$ws.Range("L7").Value = 10.37280988393497

$ws.Range("B8").Value = 14.59862207325103
$ws.Range("C8").Value = 4.482851133118324
$ws.Range("D8").Value = 8.574865400116613
$ws.Range("E8").Value = 10.07348778004926
$ws.Range("F8").Value = 43.27149609731136
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 35.32403676911836
$ws.Range("K8").Value = 12.89263647667349
$ws.Range("L8").Value = 10.37935176427846

$ws.Range("B9").Value = 14.88705600075644
$ws.Range("C9").Value = 5.026492430600872
$ws.Range("D9").Value = 8.594254565786672
$ws.Range("E9").Value = 10.06505484436621
$ws.Range("F9").Value = 44.16303208732247
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 35.76469086133142
$ws.Range("K9").Value = 13.09642447285637
$ws.Range("L9").Value = 10.40786274426256

$ws.Range("B10").Value = 15.11897000663642
$ws.Range("C10").Value = 5.38734057311506
$ws.Range("D10").Value = 8.607662949769662
$ws.Range("E10").Value = 10.06598543810546
$ws.Range("F10").Value = 44.83832896536155
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 36.10499053265374
$ws.Range("K10").Value = 13.26357003738432
$ws.Range("L10").Value = 10.43827860457699

$ws.Range("B11").Value = 15.22837739660829
$ws.Range("C11").Value = 5.543005406550329
$ws.Range("D11").Value = 8.613585109635128
$ws.Range("E11").Value = 10.06795009638252
$ws.Range("F11").Value = 45.14910073752421
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 36.26311309547312
$ws.Range("K11").Value = 13.34311096286092
$ws.Range("L11").Value = 10.45414622569835

$ws.Range("B12").Value = 15.27033045246928
$ws.Range("C12").Value = 5.600724470005993
$ws.Range("D12").Value = 8.61580247604468
$ws.Range("E12").Value = 10.06891499163133
$ws.Range("F12").Value = 45.26722323204387
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 36.3234415667719
$ws.Range("K12").Value = 13.3737098151865
$ws.Range("L12").Value = 10.46044431305005

$ws.Range("B13").Value = 15.26127249077541
$ws.Range("C13").Value = 5.588348340289196
$ws.Range("D13").Value = 8.61532604406664
$ws.Range("E13").Value = 10.06869737113984
$ws.Range("F13").Value = 45.2417650601874
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 36.31042912574069
$ws.Range("K13").Value = 13.36709894507569
$ws.Range("L13").Value = 10.45907508433295

$ws.Range("B14").Value = 15.23181864872218
$ws.Range("C14").Value = 5.547778630995757
$ws.Range("D14").Value = 8.613768037700064
$ws.Range("E14").Value = 10.06802505576856
$ws.Range("F14").Value = 45.15881035606589
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 36.26806747242107
$ws.Range("K14").Value = 13.34561891073667
$ws.Range("L14").Value = 10.45465859054433

$ws.Range("B15").Value = 15.21384423198973
$ws.Range("C15").Value = 5.522768429791072
$ws.Range("D15").Value = 8.612810436764383
$ws.Range("E15").Value = 10.06764199067719
$ws.Range("F15").Value = 45.10805328721016
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 36.24217765986665
$ws.Range("K15").Value = 13.33252329740155
$ws.Range("L15").Value = 10.45199095818194

$ws.Range("B16").Value = 15.11189529780639
$ws.Range("C16").Value = 5.376996046592769
$ws.Range("D16").Value = 8.607272377367197
$ws.Range("E16").Value = 10.06588799932018
$ws.Range("F16").Value = 44.81808492476303
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 36.0947213995295
$ws.Range("K16").Value = 13.25844021986119
$ws.Range("L16").Value = 10.43728226596283

$ws.Range("B17").Value = 15.05032642878158
$ws.Range("C17").Value = 5.285389715672673
$ws.Range("D17").Value = 8.60382972625607
$ws.Range("E17").Value = 10.06520627333459
$ws.Range("F17").Value = 44.64106200681069
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 36.00509327827893
$ws.Range("K17").Value = 13.21387256300096
$ws.Range("L17").Value = 10.42877736986075

$ws.Range("B18").Value = 15.01528395070773
$ws.Range("C18").Value = 5.231901645929173
$ws.Range("D18").Value = 8.601832891910275
$ws.Range("E18").Value = 10.06495928893036
$ws.Range("F18").Value = 44.53958418092946
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 35.95385554367743
$ws.Range("K18").Value = 13.18856998981829
$ws.Range("L18").Value = 10.42407695578406

$ws.Range("B19").Value = 15.00348394855201
$ws.Range("C19").Value = 5.213654620134404
$ws.Range("D19").Value = 8.601153919125053
$ws.Range("E19").Value = 10.06490061010264
$ws.Range("F19").Value = 44.50528636231977
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 35.93656204238705
$ws.Range("K19").Value = 13.18006068579751
$ws.Range("L19").Value = 10.42251842409047

$ws.Range("B20").Value = 15.05684249563256
$ws.Range("C20").Value = 5.295224071800052
$ws.Range("D20").Value = 8.604197929994106
$ws.Range("E20").Value = 10.06526382910161
$ws.Range("F20").Value = 44.65987163447094
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 36.01460203600157
$ws.Range("K20").Value = 13.2185827334672
$ws.Range("L20").Value = 10.42966294265956

$ws.Range("B21").Value = 15.24045608285906
$ws.Range("C21").Value = 5.559728309406191
$ws.Range("D21").Value = 8.614226344374893
$ws.Range("E21").Value = 10.06821654156514
$ws.Range("F21").Value = 45.18316481762339
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 36.28049809384135
$ws.Range("K21").Value = 13.35191534222234
$ws.Range("L21").Value = 10.4559479918032

$ws.Range("B22").Value = 15.3634870029962
$ws.Range("C22").Value = 5.725438409135882
$ws.Range("D22").Value = 8.620633578023391
$ws.Range("E22").Value = 10.0714336066257
$ws.Range("F22").Value = 45.52769683616779
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 36.45689187034045
$ws.Range("K22").Value = 13.44183203980072
$ws.Range("L22").Value = 10.47481195257578

$ws.Range("B23").Value = 15.29755911430304
$ws.Range("C23").Value = 5.637652738418064
$ws.Range("D23").Value = 8.617227267956668
$ws.Range("E23").Value = 10.06959907422488
$ws.Range("F23").Value = 45.34360681924257
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 36.36251679958548
$ws.Range("K23").Value = 13.393596502061
$ws.Range("L23").Value = 10.46459068685415

$ws.Range("B24").Value = 15.05389547528802
$ws.Range("C24").Value = 5.290780520355922
$ws.Range("D24").Value = 8.604031519944481
$ws.Range("E24").Value = 10.06523735659935
$ws.Range("F24").Value = 44.65136687718056
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 36.01030221919742
$ws.Range("K24").Value = 13.21645226680448
$ws.Range("L24").Value = 10.42926198578016

$ws.Range("B25").Value = 14.80536995359446
$ws.Range("C25").Value = 4.886144561423161
$ws.Range("D25").Value = 8.589157607797462
$ws.Range("E25").Value = 10.06608289690145
$ws.Range("F25").Value = 43.91800344213559
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 35.64248276612115
$ws.Range("K25").Value = 13.03814341453413
$ws.Range("L25").Value = 10.39847915149024
